# Add "270 Packets" run data into columns I:K (rows 3-51),
# add a new "250 Packets" header block (N1 + M2:O2 headers),
# matching the author's commit "added 270 packet run".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New group header in row 1: "250 Packets" over columns M:O ---
$ws.Range("N1").Value = "250 Packets"

# --- New group sub-headers in row 2: Acc / Loss / Time for columns M:O ---
$ws.Range("M2").Value = "Acc"
$ws.Range("N2").Value = "Loss"
$ws.Range("O2").Value = "Time"

# --- New "270 Packets" data for columns I:K, rows 3-51 ---
$data = @(
  @(89.674389362335205, 0.244528188925136,   82.455614089965806),
  @(89.342379570007296, 0.25267978321172402, 75.036825656890798),
  @(88.751870393752995, 0.25918662667721099, 75.991146802902193),
  @(89.909172058105398, 0.25369177561717199, 76.640452146530095),
  @(89.185857772827106, 0.24768675051806699, 76.136749982833805),
  @(88.241988420486393, 0.26343331673686798, 75.151413202285696),
  @(88.903647661209106, 0.25511473001401602, 74.383100032806396),
  @(88.804042339324894, 0.25308474167839601, 74.369613409042302),
  @(88.585859537124605, 0.25379235884448298, 74.491253614425602),
  @(89.795339107513399, 0.23951616773269299, 75.908175230026202),
  @(89.937627315521198, 0.24038724478290899, 75.984380960464406),
  @(88.585859537124605, 0.25078049290404097, 75.073441743850694),
  @(90.295726060867295, 0.23855006397332301, 75.772522926330495),
  @(89.254629611968994, 0.261448279177599,   75.891891717910696),
  @(88.747125864028902, 0.25448242167212898, 74.464914083480807),
  @(89.418268203735295, 0.24392899600639401, 74.943201541900606),
  @(89.304435253143296, 0.24944798000179599, 76.312826156616197),
  @(89.444351196289006, 0.24763306652123901, 76.144528627395601),
  @(88.832497596740694, 0.25881855143789601, 75.318665742874103),
  @(89.577156305313096, 0.242277082595513,   75.418741703033405),
  @(88.9131307601928,   0.25252829890109002, 75.688076972961397),
  @(89.311546087265,    0.25169569942632403, 76.163012027740393),
  @(88.972419500350895, 0.25440847281804102, 75.628634691238403),
  @(88.889414072036701, 0.25551032082455299, 75.975719213485704),
  @(89.192968606948796, 0.25626749786626601, 76.8635346889495),
  @(88.972419500350895, 0.25850718530676298, 76.635301589965806),
  @(89.294946193694997, 0.25150514153655101, 75.884287595748901),
  @(89.1716241836547,   0.25054168452265502, 74.442911386489797),
  @(88.863331079483004, 0.27117592397217899, 75.919968605041504),
  @(89.577156305313096, 0.244464064881007,   76.845924854278493),
  @(89.197713136672903, 0.26747144625227998, 76.167885541915894),
  @(88.732898235321002, 0.264799719276476,   76.130175828933702),
  @(88.699692487716604, 0.258527421027829,   76.171753644943195),
  @(89.534467458724905, 0.25453679464543799, 76.952017784118596),
  @(89.257001876830998, 0.24804810150953499, 76.221291780471802),
  @(89.577156305313096, 0.25602399755423799, 75.4796013832092),
  @(89.005619287490802, 0.258948004286644,   75.752249002456594),
  @(89.226174354553194, 0.25585194801063998, 76.804516315460205),
  @(89.3637180328369,   0.23998273261691599, 76.135252714157104),
  @(89.570045471191406, 0.244280452884269,   76.395229101181002),
  @(88.929730653762803, 0.25502020090282601, 75.890778541564899),
  @(89.067280292510901, 0.25527764901182698, 75.688072443008394),
  @(88.955819606781006, 0.25693635091425499, 76.632386922836304),
  @(89.669644832610999, 0.242593294158473,   75.193731069564805),
  @(89.396923780441199, 0.24983933237144601, 76.821705341339097),
  @(88.8846755027771,   0.25353088884362601, 76.058594942092896),
  @(88.955819606781006, 0.25039767996927298, 76.064940690994206),
  @(89.396923780441199, 0.25604076211240601, 75.976545572280799),
  @(89.439606666564899, 0.25780443963407501, 75.536620140075598)
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $startRow + $i
  $ws.Cells.Item($row, 9).Value  = $data[$i][0]   # column I
  $ws.Cells.Item($row, 10).Value = $data[$i][1]   # column J
  $ws.Cells.Item($row, 11).Value = $data[$i][2]   # column K
}

# --- Match the recorded selection after the edit ---
$ws.Range("M3").Select() | Out-Null
